$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 478.875  # H6: 1500 -> 478.875
$ws.Cells.Item(6, 9).Value = 478.875  # I6: 1500 -> 478.875
$ws.Cells.Item(6, 11).Value = 1436.625  # K6: 4500 -> 1436.625
$ws.Cells.Item(6, 13).Value = -1324.625  # M6: -4388 -> -1324.625

# Sheet ALC, Row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 55.666668  # H9: 44.25 -> 55.666668
$ws.Cells.Item(9, 9).Value = 20  # I9: 15 -> 20
$ws.Cells.Item(9, 11).Value = 20  # K9: 15 -> 20
$ws.Cells.Item(9, 13).Value = 149  # M9: 154 -> 149

# Sheet ALC, Row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 189.8  # H12: 127.25 -> 189.8
$ws.Cells.Item(12, 9).Value = 187.25  # I12: 127.5 -> 187.25
$ws.Cells.Item(12, 10).Value = 200  # J12: 127 -> 200
$ws.Cells.Item(12, 11).Value = 187.25  # K12: 127.5 -> 187.25
$ws.Cells.Item(12, 12).Value = 200  # L12: 127 -> 200
$ws.Cells.Item(12, 13).Value = -17.25  # M12: 42.5 -> -17.25
$ws.Cells.Item(12, 14).Value = -540  # N12: -467 -> -540

# Sheet ALC, Row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 0  # H21: 21250 -> 0
$ws.Cells.Item(21, 9).Value = 0  # I21: 25000 -> 0
$ws.Cells.Item(21, 10).Value = 0  # J21: 20000 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 25000 -> 0
$ws.Cells.Item(21, 12).Value = 0  # L21: 20000 -> 0
$ws.Cells.Item(21, 13).ClearContents()  # M21: -24532 -> (removed)
$ws.Cells.Item(21, 14).ClearContents()  # N21: -20936 -> (removed)

# Sheet ALC, Row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(23, 8).Value = 0  # H23: 21250 -> 0
$ws.Cells.Item(23, 9).Value = 0  # I23: 25000 -> 0
$ws.Cells.Item(23, 10).Value = 0  # J23: 20000 -> 0
$ws.Cells.Item(23, 11).Value = 0  # K23: 25000 -> 0
$ws.Cells.Item(23, 12).Value = 0  # L23: 20000 -> 0
$ws.Cells.Item(23, 13).ClearContents()  # M23: -24766 -> (removed)
$ws.Cells.Item(23, 14).ClearContents()  # N23: -20468 -> (removed)

# Sheet ALC, Row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 275  # H29: 1250 -> 275
$ws.Cells.Item(29, 9).Value = 275  # I29: 1000 -> 275
$ws.Cells.Item(29, 10).Value = 0  # J29: 3000 -> 0
$ws.Cells.Item(29, 11).Value = 825  # K29: 3000 -> 825
$ws.Cells.Item(29, 12).Value = 0  # L29: 9000 -> 0
$ws.Cells.Item(29, 13).Value = -544  # M29: -2719 -> -544
$ws.Cells.Item(29, 14).ClearContents()  # N29: -9562 -> (removed)

# Sheet ALC, Row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 1771.75  # H100: 1485.4445 -> 1771.75
$ws.Cells.Item(100, 9).Value = 1500.4  # I100: 1006.3571 -> 1500.4
$ws.Cells.Item(100, 10).Value = 1965.5714  # J100: 2001.3846 -> 1965.5714
$ws.Cells.Item(100, 11).Value = 1500.4  # K100: 1006.3571 -> 1500.4
$ws.Cells.Item(100, 12).Value = 1965.5714  # L100: 2001.3846 -> 1965.5714
$ws.Cells.Item(100, 13).Value = -959.4000000000001  # M100: -465.3570999999999 -> -959.4000000000001
$ws.Cells.Item(100, 14).Value = -3047.5714  # N100: -3083.3846 -> -3047.5714

# Sheet ALC, Row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2441.1177  # H116: 2559.9333 -> 2441.1177
$ws.Cells.Item(116, 9).Value = 2224.9167  # I116: 2318.0908 -> 2224.9167
$ws.Cells.Item(116, 10).Value = 2960  # J116: 3225 -> 2960
$ws.Cells.Item(116, 11).Value = 2224.9167  # K116: 2318.0908 -> 2224.9167
$ws.Cells.Item(116, 12).Value = 2960  # L116: 3225 -> 2960
$ws.Cells.Item(116, 13).Value = 1217.0833  # M116: 1123.9092 -> 1217.0833
$ws.Cells.Item(116, 14).Value = -9844  # N116: -10109 -> -9844

# Sheet ALC, Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 6414001.5  # H132: 4466903 -> 6414001.5
$ws.Cells.Item(132, 9).Value = 3886  # I132: 2533.34 -> 3886
$ws.Cells.Item(132, 10).Value = 35717384  # J132: 41669984 -> 35717384
$ws.Cells.Item(132, 11).Value = 11658  # K132: 7600.02 -> 11658
$ws.Cells.Item(132, 12).Value = 107152152  # L132: 125009952 -> 107152152
$ws.Cells.Item(132, 13).Value = -9128  # M132: -5070.02 -> -9128
$ws.Cells.Item(132, 14).Value = -107157212  # N132: -125015012 -> -107157212

# Sheet ARM, Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 18740.77  # H32: 16134.767 -> 18740.77
$ws.Cells.Item(32, 9).Value = 18593.83  # I32: 15719.754 -> 18593.83
$ws.Cells.Item(32, 11).Value = 18593.83  # K32: 15719.754 -> 18593.83
$ws.Cells.Item(32, 13).Value = -18306.83  # M32: -15432.754 -> -18306.83

# Sheet ARM, Row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1123.3334  # H102: 1500 -> 1123.3334
$ws.Cells.Item(102, 9).Value = 1113.75  # I102: 1500 -> 1113.75
$ws.Cells.Item(102, 10).Value = 1200  # J102: 0 -> 1200
$ws.Cells.Item(102, 11).Value = 1113.75  # K102: 1500 -> 1113.75
$ws.Cells.Item(102, 12).Value = 1200  # L102: 0 -> 1200
$ws.Cells.Item(102, 13).Value = 508.25  # M102: 122 -> 508.25
$ws.Cells.Item(102, 14).Value = -4444  # N102: None -> -4444

# Sheet CRP, Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4529.476  # H31: 3638.7637 -> 4529.476
$ws.Cells.Item(31, 9).Value = 3237.375  # I31: 2170.4 -> 3237.375
$ws.Cells.Item(31, 10).Value = 5324.615  # J31: 5400.8 -> 5324.615
$ws.Cells.Item(31, 11).Value = 3237.375  # K31: 2170.4 -> 3237.375
$ws.Cells.Item(31, 12).Value = 5324.615  # L31: 5400.8 -> 5324.615
$ws.Cells.Item(31, 13).Value = -2942.375  # M31: -1875.4 -> -2942.375
$ws.Cells.Item(31, 14).Value = -5914.615  # N31: -5990.8 -> -5914.615

# Sheet CRP, Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 4529.476  # H34: 3638.7637 -> 4529.476
$ws.Cells.Item(34, 9).Value = 3237.375  # I34: 2170.4 -> 3237.375
$ws.Cells.Item(34, 10).Value = 5324.615  # J34: 5400.8 -> 5324.615
$ws.Cells.Item(34, 11).Value = 3237.375  # K34: 2170.4 -> 3237.375
$ws.Cells.Item(34, 12).Value = 5324.615  # L34: 5400.8 -> 5324.615
$ws.Cells.Item(34, 13).Value = -3035.375  # M34: -1968.4 -> -3035.375
$ws.Cells.Item(34, 14).Value = -5728.615  # N34: -5804.8 -> -5728.615

# Sheet CUL, Row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 32259434  # H113: 29413066 -> 32259434
$ws.Cells.Item(113, 9).Value = 76924380  # I113: 66667884 -> 76924380
$ws.Cells.Item(113, 10).Value = 1411.1111  # J113: 1368.421 -> 1411.1111
$ws.Cells.Item(113, 11).Value = 230773140  # K113: 200003652 -> 230773140
$ws.Cells.Item(113, 12).Value = 4233.3333  # L113: 4105.263 -> 4233.3333
$ws.Cells.Item(113, 13).Value = -230770970  # M113: -200001482 -> -230770970
$ws.Cells.Item(113, 14).Value = -8573.3333  # N113: -8445.262999999999 -> -8573.3333

# Sheet CUL, Row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 1674.381  # H122: 1604.25 -> 1674.381
$ws.Cells.Item(122, 10).Value = 1207.9  # J122: 1195.5652 -> 1207.9
$ws.Cells.Item(122, 12).Value = 10871.1  # L122: 10760.0868 -> 10871.1
$ws.Cells.Item(122, 14).Value = -15771.1  # N122: -15660.0868 -> -15771.1

# Sheet CUL, Row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(125, 8).Value = 4086.9355  # H125: 4046.7188 -> 4086.9355
$ws.Cells.Item(125, 10).Value = 4569.375  # J125: 4498.6 -> 4569.375
$ws.Cells.Item(125, 12).Value = 13708.125  # L125: 13495.8 -> 13708.125
$ws.Cells.Item(125, 14).Value = -23548.125  # N125: -23335.8 -> -23548.125

# Sheet GSM, Row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 75.75  # H2: 105.42857 -> 75.75
$ws.Cells.Item(2, 9).Value = 75.7  # I2: 109 -> 75.7
$ws.Cells.Item(2, 10).Value = 76  # J2: 100.666664 -> 76
$ws.Cells.Item(2, 11).Value = 75.7  # K2: 109 -> 75.7
$ws.Cells.Item(2, 12).Value = 76  # L2: 100.666664 -> 76
$ws.Cells.Item(2, 13).Value = 37.3  # M2: 4 -> 37.3
$ws.Cells.Item(2, 14).Value = -302  # N2: -326.666664 -> -302

# Sheet GSM, Row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 1375  # H43: 6937.5 -> 1375
$ws.Cells.Item(43, 9).Value = 1375  # I43: 1166.6666 -> 1375
$ws.Cells.Item(43, 10).Value = 0  # J43: 10400 -> 0
$ws.Cells.Item(43, 11).Value = 1375  # K43: 1166.6666 -> 1375
$ws.Cells.Item(43, 12).Value = 0  # L43: 10400 -> 0
$ws.Cells.Item(43, 13).Value = -1224  # M43: -1015.6666 -> -1224
$ws.Cells.Item(43, 14).ClearContents()  # N43: -10702 -> (removed)

# Sheet GSM, Row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 15198.333  # H46: 18475.1 -> 15198.333
$ws.Cells.Item(46, 9).Value = 1299.3334  # I46: 2190 -> 1299.3334
$ws.Cells.Item(46, 10).Value = 22147.834  # J46: 22546.375 -> 22147.834
$ws.Cells.Item(46, 11).Value = 1299.3334  # K46: 2190 -> 1299.3334
$ws.Cells.Item(46, 12).Value = 22147.834  # L46: 22546.375 -> 22147.834
$ws.Cells.Item(46, 13).Value = -1143.3334  # M46: -2034 -> -1143.3334
$ws.Cells.Item(46, 14).Value = -22459.834  # N46: -22858.375 -> -22459.834

# Sheet GSM, Row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 9596.556  # H57: 6400.8335 -> 9596.556
$ws.Cells.Item(57, 9).Value = 4479.8  # I57: 4481 -> 4479.8
$ws.Cells.Item(57, 10).Value = 15992.5  # J57: 16000 -> 15992.5
$ws.Cells.Item(57, 11).Value = 4479.8  # K57: 4481 -> 4479.8
$ws.Cells.Item(57, 12).Value = 15992.5  # L57: 16000 -> 15992.5
$ws.Cells.Item(57, 13).Value = -3659.8  # M57: -3661 -> -3659.8
$ws.Cells.Item(57, 14).Value = -17632.5  # N57: -17640 -> -17632.5

# Sheet GSM, Row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2914.4285  # H80: 2769.05 -> 2914.4285
$ws.Cells.Item(80, 9).Value = 2743.5715  # I80: 2633.5833 -> 2743.5715
$ws.Cells.Item(80, 10).Value = 3085.2856  # J80: 2972.25 -> 3085.2856
$ws.Cells.Item(80, 11).Value = 2743.5715  # K80: 2633.5833 -> 2743.5715
$ws.Cells.Item(80, 12).Value = 3085.2856  # L80: 2972.25 -> 3085.2856
$ws.Cells.Item(80, 13).Value = -1745.5715  # M80: -1635.5833 -> -1745.5715
$ws.Cells.Item(80, 14).Value = -5081.2856  # N80: -4968.25 -> -5081.2856

# Sheet GSM, Row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2914.4285  # H83: 2769.05 -> 2914.4285
$ws.Cells.Item(83, 9).Value = 2743.5715  # I83: 2633.5833 -> 2743.5715
$ws.Cells.Item(83, 10).Value = 3085.2856  # J83: 2972.25 -> 3085.2856
$ws.Cells.Item(83, 11).Value = 13717.8575  # K83: 13167.9165 -> 13717.8575
$ws.Cells.Item(83, 12).Value = 15426.428  # L83: 14861.25 -> 15426.428
$ws.Cells.Item(83, 13).Value = -8725.8575  # M83: -8175.916499999999 -> -8725.8575
$ws.Cells.Item(83, 14).Value = -25410.428  # N83: -24845.25 -> -25410.428

# Sheet GSM, Row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1864.0769  # H97: 1980.2727 -> 1864.0769
$ws.Cells.Item(97, 9).Value = 1350  # I97: 1433.3334 -> 1350
$ws.Cells.Item(97, 11).Value = 1350  # K97: 1433.3334 -> 1350
$ws.Cells.Item(97, 13).Value = -854  # M97: -937.3334 -> -854

# Sheet GSM, Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 6895.0415  # H132: 7020.0415 -> 6895.0415
$ws.Cells.Item(132, 10).Value = 3998.75  # J132: 4748.75 -> 3998.75
$ws.Cells.Item(132, 12).Value = 11996.25  # L132: 14246.25 -> 11996.25
$ws.Cells.Item(132, 14).Value = -17056.25  # N132: -19306.25 -> -17056.25

# Sheet LTW, Row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8966.522999999999  # H7: 9190.809999999999 -> 8966.522999999999
$ws.Cells.Item(7, 9).Value = 6163.5454  # I7: 6700.9 -> 6163.5454
$ws.Cells.Item(7, 10).Value = 12049.8  # J7: 11454.363 -> 12049.8
$ws.Cells.Item(7, 11).Value = 6163.5454  # K7: 6700.9 -> 6163.5454
$ws.Cells.Item(7, 12).Value = 12049.8  # L7: 11454.363 -> 12049.8
$ws.Cells.Item(7, 13).Value = -6051.5454  # M7: -6588.9 -> -6051.5454
$ws.Cells.Item(7, 14).Value = -12273.8  # N7: -11678.363 -> -12273.8

# Sheet LTW, Row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3093.2727  # H40: 2790.1538 -> 3093.2727
$ws.Cells.Item(40, 9).Value = 8740.429  # I40: 5965.5454 -> 8740.429
$ws.Cells.Item(40, 10).Value = 1572.8846  # J40: 1542.6786 -> 1572.8846
$ws.Cells.Item(40, 11).Value = 8740.429  # K40: 5965.5454 -> 8740.429
$ws.Cells.Item(40, 12).Value = 1572.8846  # L40: 1542.6786 -> 1572.8846
$ws.Cells.Item(40, 13).Value = -8604.429  # M40: -5829.5454 -> -8604.429
$ws.Cells.Item(40, 14).Value = -1844.8846  # N40: -1814.6786 -> -1844.8846

# Sheet LTW, Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 9176.352999999999  # H122: 15077.556 -> 9176.352999999999
$ws.Cells.Item(122, 9).Value = 10083.333  # I122: 50000 -> 10083.333
$ws.Cells.Item(122, 10).Value = 8681.637000000001  # J122: 10712.25 -> 8681.637000000001
$ws.Cells.Item(122, 11).Value = 30249.999  # K122: 150000 -> 30249.999
$ws.Cells.Item(122, 12).Value = 26044.911  # L122: 32136.75 -> 26044.911
$ws.Cells.Item(122, 13).Value = -27799.999  # M122: -147550 -> -27799.999
$ws.Cells.Item(122, 14).Value = -30944.911  # N122: -37036.75 -> -30944.911

# Sheet LTW, Row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 8966.522999999999  # H126: 9190.809999999999 -> 8966.522999999999
$ws.Cells.Item(126, 9).Value = 6163.5454  # I126: 6700.9 -> 6163.5454
$ws.Cells.Item(126, 10).Value = 12049.8  # J126: 11454.363 -> 12049.8
$ws.Cells.Item(126, 11).Value = 18490.6362  # K126: 20102.7 -> 18490.6362
$ws.Cells.Item(126, 12).Value = 36149.39999999999  # L126: 34363.089 -> 36149.39999999999
$ws.Cells.Item(126, 13).Value = -16020.6362  # M126: -17632.7 -> -16020.6362
$ws.Cells.Item(126, 14).Value = -41089.39999999999  # N126: -39303.089 -> -41089.39999999999

# Sheet WVR, Row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 27420  # H54: 21365 -> 27420
$ws.Cells.Item(54, 9).Value = 19690  # I54: 7220 -> 19690
$ws.Cells.Item(54, 10).Value = 29996.666  # J54: 27427.143 -> 29996.666
$ws.Cells.Item(54, 11).Value = 19690  # K54: 7220 -> 19690
$ws.Cells.Item(54, 12).Value = 29996.666  # L54: 27427.143 -> 29996.666
$ws.Cells.Item(54, 13).Value = -19170  # M54: -6700 -> -19170
$ws.Cells.Item(54, 14).Value = -31036.666  # N54: -28467.143 -> -31036.666

# Sheet WVR, Row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5469.231  # H62: 5371.875 -> 5469.231
$ws.Cells.Item(62, 10).Value = 5314.2856  # J62: 5205 -> 5314.2856
$ws.Cells.Item(62, 12).Value = 5314.2856  # L62: 5205 -> 5314.2856
$ws.Cells.Item(62, 14).Value = -6562.2856  # N62: -6453 -> -6562.2856

# Sheet WVR, Row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 5469.231  # H65: 5371.875 -> 5469.231
$ws.Cells.Item(65, 10).Value = 5314.2856  # J65: 5205 -> 5314.2856
$ws.Cells.Item(65, 12).Value = 26571.428  # L65: 26025 -> 26571.428
$ws.Cells.Item(65, 14).Value = -32811.428  # N65: -32265 -> -32811.428
